# Update Chris Gayle's innings-by-innings batting stats (runs, balls, fours,
# sixes) so that rows 2,4,5,6,7,8 take on the values described in the commit
# diff. Row 3 (24,21,1,2) is left untouched.
#
# The source values are stored as text (t="str"/shared-string) even though
# they look numeric, so we force each target cell to Text format before
# writing, otherwise Excel would silently reinterpret "53" etc. as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @("53", "45", "1", "5")
    4 = @("29", "13", "3", "2")
    5 = @("51", "29", "2", "5")
    6 = @("12", "19", "0", "0")
    7 = @("99", "63", "6", "8")
    8 = @("20", "20", "2", "1")
}

$columns = @("C", "D", "E", "F")

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cell = $ws.Range("$($columns[$i])$row")
        $cell.NumberFormat = "@"
        $cell.Value = $values[$i]
        # Drop the explicit text number-format override again so the cell
        # keeps using the default style, just like every other cell in the
        # sheet (only the stored value should differ from the original).
        $cell.Style = "Normal"
    }
}
